# feat: add 2022-Q1 data
#
# The existing "总计" (totals) sheet is repurposed into the new "2022-Q1"
# per-fund holdings sheet (it keeps sheetId/rId 6), while a fresh "总计"
# sheet is appended at the end (sheetId/rId 7) holding the refreshed totals
# table -- the old rows shifted down by one, with a new first data row for
# 2022-Q1.

$wb = $excel.ActiveWorkbook

$q4 = $wb.Worksheets.Item("2021-Q4")
$total = $wb.Worksheets.Item("总计")

# Data captured from the previous "总计" table before it gets overwritten.
$oldDates  = @("2021-Q4", "2021-Q3", "2021-Q2", "2021-Q1", "2020-Q4")
$oldCounts = @(3, 4, 4, 6, 2)
$oldValues = @(0.89, 1.24, 1.28, 1.46, 1.21)

# --- Step 1: copy "总计" to the end of the workbook first, *before* it is
#     renamed -- this gives the brand-new totals sheet the right sheetPr /
#     pageMargins / column layout "for free", since it is cloned from a
#     sheet that already has that exact 4-column totals shape. ---
$total.Copy($null, $total)
$newTotal = $wb.Worksheets.Item($wb.Worksheets.Count)

# --- Step 2: rename the original "总计" sheet to "2022-Q1" and the copy
#     back to "总计", so sheetId/rId land as 6 and 7 respectively. ---
$total.Name = "2022-Q1"
$newTotal.Name = "总计"
$newQ1 = $total

# --- Step 3: replace "2022-Q1" contents with the new per-fund holdings
#     table (same shape/style as the "2021-Q4" sheet). ---
$newQ1.Cells.Clear()
$q4.Range("A1:H4").Copy($newQ1.Range("A1"))

$newQ1.Range("B1").Value = "基金代码"
$newQ1.Range("C1").Value = "基金名称"
$newQ1.Range("D1").Value = "基金规模"
$newQ1.Range("E1").Value = "股票总仓位"
$newQ1.Range("F1").Value = "仓位占比"
$newQ1.Range("G1").Value = "持有市值(亿元)"
$newQ1.Range("H1").Value = "仓位排名"

$newQ1.Range("B2").Value = "'000480"
$newQ1.Range("C2").Value = "东方红新动力灵活配置混合"
$newQ1.Range("D2").Value = "'15.38"
$newQ1.Range("E2").Value = "'72.90"
$newQ1.Range("F2").Value = "'3.92"
$newQ1.Range("G2").Value = "'0.6029"
$newQ1.Range("H2").Value = 5

$newQ1.Range("B3").Value = "'001564"
$newQ1.Range("C3").Value = "东方红京东大数据灵活配置混合"
$newQ1.Range("D3").Value = "'11.40"
$newQ1.Range("E3").Value = "'69.58"
$newQ1.Range("F3").Value = "'2.63"
$newQ1.Range("G3").Value = "'0.2998"
$newQ1.Range("H3").Value = 8

$newQ1.Range("B4").Value = "'002367"
$newQ1.Range("C4").Value = "国联安安稳灵活配置混合"
$newQ1.Range("D4").Value = "'2.32"
$newQ1.Range("E4").Value = "'33.99"
$newQ1.Range("F4").Value = "'2.12"
$newQ1.Range("G4").Value = "'0.0492"
$newQ1.Range("H4").Value = 5

# --- Step 4: refresh the new "总计" sheet -- insert a 2022-Q1 row at the
#     top and push the previously-captured rows down by one. ---
$newTotal.Range("B1").Value = "日期"
$newTotal.Range("C1").Value = "持有数量(只)"
$newTotal.Range("D1").Value = "持有市值(亿元)"

$newTotal.Range("A2").Value = 0
$newTotal.Range("B2").Value = "2022-Q1"
$newTotal.Range("C2").Value = 3
$newTotal.Range("D2").Value = 0.95

for ($i = 0; $i -lt $oldDates.Length; $i++) {
    $row = $i + 3
    if ($row -gt 6) {
        # Row 7 doesn't exist on the copied sheet yet (the original "总计"
        # only had 6 rows) -- clone the A-column style from the row above
        # so the new row-number cell keeps the same "s=2" formatting.
        $newTotal.Cells.Item($row - 1, 1).Copy($newTotal.Cells.Item($row, 1))
    }
    $newTotal.Cells.Item($row, 1).Value = $i + 1
    $newTotal.Cells.Item($row, 2).Value = $oldDates[$i]
    $newTotal.Cells.Item($row, 3).Value = $oldCounts[$i]
    $newTotal.Cells.Item($row, 4).Value = $oldValues[$i]
}
